$d = $word.ActiveDocument

# Move to the very end of the document and insert one more paragraph break,
# matching the last (empty) paragraph that already exists there.
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd = 0
$endRange.InsertParagraphAfter()
